# "updated list of metrics"
#
# Row 20 of Sheet1 documents the metrics_kde() function. The metric-name
# list and its notes/description were revised:
#   - Metrics name (A20): the old 3-value list of outputs is replaced with
#     the new, expanded list describing multiple peaks.
#   - Description (B20): the old long paragraph explaining KDE is replaced
#     with a short note about how the new version differs from Fusion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "kde_peaks_count, kde_peak1_elev,  kde_peak2_elev, …, kde_peak1_value, kde_peak2_value, …, kde_peak1_diff, kde_peak2_diff, …"
$ws.Range("B20").Value = "Based on similar metric available in Fusion (see references), with significant differences in the list of output statistics as well as the default bandwidth used when estimating kernel density."

# Update the view state to match: scrolled so row 15 is at the top, with
# B22 selected as the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("B22").Select()
